# Weekly NYPD CompStat update: refresh report header dates/volume number
# and overwrite the crime-stat tables for rows 15-31 with newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header (Volume/Number banner and reporting week dates) ---
# These cells are rich-text shared strings with uniform run formatting, so
# re-writing the whole cell text preserves the visual formatting.
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# --- Precinct crime-stat table: plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 15
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -25
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 87.5
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 23.529411764705
$ws.Range("I16").Value = 162
$ws.Range("J16").Value = 187
$ws.Range("K16").Value = -13.368983957219
$ws.Range("L16").Value = -18.181818181818
$ws.Range("M16").Value = 8
$ws.Range("N16").Value = -79.015544041450
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 4.545454545454
$ws.Range("I17").Value = 227
$ws.Range("J17").Value = 202
$ws.Range("K17").Value = 12.376237623762
$ws.Range("L17").Value = 20.744680851063
$ws.Range("M17").Value = 136.458333333333
$ws.Range("N17").Value = -4.219409282700
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -10.526315789473
$ws.Range("I18").Value = 171
$ws.Range("J18").Value = 199
$ws.Range("K18").Value = -14.070351758794
$ws.Range("L18").Value = 6.211180124223
$ws.Range("M18").Value = -8.064516129032
$ws.Range("N18").Value = -85.434412265758
$ws.Range("C19").Value = 22
$ws.Range("E19").Value = 46.666666666666
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 42.307692307692
$ws.Range("I19").Value = 615
$ws.Range("J19").Value = 577
$ws.Range("K19").Value = 6.585788561525
$ws.Range("L19").Value = 6.401384083044
$ws.Range("M19").Value = 78.260869565217
$ws.Range("N19").Value = -13.744740532959
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 25
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 185
$ws.Range("J20").Value = 202
$ws.Range("K20").Value = -8.415841584158
$ws.Range("L20").Value = -18.141592920354
$ws.Range("M20").Value = 12.804878048780
$ws.Range("N20").Value = -88.276299112801
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 2.777777777777
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = 18.382352941176
$ws.Range("I21").Value = 1377
$ws.Range("J21").Value = 1383
$ws.Range("K21").Value = -0.433839479392
$ws.Range("L21").Value = 0.437636761487
$ws.Range("M21").Value = 44.188481675392
$ws.Range("N21").Value = -69.331848552338
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 46
$ws.Range("K22").Value = 13.043478260869
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = 26.829268292682
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -6.896551724137
$ws.Range("F24").Value = 123
$ws.Range("H24").Value = -15.172413793103
$ws.Range("I24").Value = 1212
$ws.Range("J24").Value = 1592
$ws.Range("K24").Value = -23.869346733668
$ws.Range("L24").Value = -21.502590673575
$ws.Range("M24").Value = 66.483516483516
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -40.909090909090
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 93
$ws.Range("H25").Value = -38.709677419354
$ws.Range("I25").Value = 642
$ws.Range("J25").Value = 1083
$ws.Range("K25").Value = -40.720221606648
$ws.Range("L25").Value = -32.985386221294
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -41.666666666666
$ws.Range("F26").Value = 56
$ws.Range("H26").Value = 47.368421052631
$ws.Range("I26").Value = 411
$ws.Range("J26").Value = 427
$ws.Range("K26").Value = -3.747072599531
$ws.Range("L26").Value = 3.007518796992
$ws.Range("M26").Value = 3.266331658291
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = -4.347826086956
$ws.Range("L27").Value = -12
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 80
$ws.Range("I28").Value = 60
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 15.384615384615
$ws.Range("L28").Value = -10.447761194029
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = -75
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = -75

# --- Cells that change "type" (blank/N-A text <-> real number) ---
# For these, we must also swap the number format so the resulting style
# index matches a genuine numeric (or text) cell, not just overwrite the value.

# C15: was text "0" (style 13) -> becomes number 1 (style 14)
$ws.Range("J15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

# C27: was text "0" (style 13) -> becomes number 1 (style 14)
$ws.Range("J15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

# C22: was number 1 (style 14) -> becomes text "0" (style 13)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# D29/G29/D30/G30: were text "0" (style 13) -> become number 2 (style 14)
$ws.Range("J15").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 2
$ws.Range("J15").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 2
$ws.Range("J15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 2
$ws.Range("J15").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 2

# E29/H29/E30/H30: were text "***.*" (style 13) -> become number -100 (style 15)
$ws.Range("L29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("L29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("L29").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("L29").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100

# D31: was number 1 (style 14) -> becomes text "0" (style 13)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("C31").Copy()
$ws.Range("D31").PasteSpecial(-4122)

# E31: was number -100 (style 15) -> becomes text "***.*" (style 13)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("M31").Copy()
$ws.Range("E31").PasteSpecial(-4122)
